$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.679513333333333
$ws.Range("H2").Value = 5.038539999999999
$ws.Range("I2").Value = 0.376631045782902
$ws.Range("J2").Value = 0.3928053077137587
$ws.Range("M2").Value = 0.1798956666666667
$ws.Range("N2").Value = 0.539687
$ws.Range("O2").Value = 0.01070918551864568
$ws.Range("P2").Value = 0.01088280728907136
$ws.Range("Q2").Value = 0.3021371707755555
$ws.Range("R2").Value = 2.71923453698
$ws.Range("S2").Value = 0.004033411741370631
$ws.Range("T2").Value = 0.004274824465973213
$ws.Range("G3").Value = 1.679513333333333
$ws.Range("H3").Value = 5.038539999999999
$ws.Range("I3").Value = 0.376631045782902
$ws.Range("J3").Value = 0.3928053077137587
$ws.Range("N3").Value = 47.39813
$ws.Range("O3").Value = 0.9405365839956962
$ws.Range("P3").Value = 0.9557849543390003
$ws.Range("Q3").Value = 26.53526377002222
$ws.Range("R3").Value = 238.8173739302
$ws.Range("S3").Value = 0.3542352772273772
$ws.Range("T3").Value = 0.3754374030973118
$ws.Range("G4").Value = 1.679513333333333
$ws.Range("H4").Value = 5.038539999999999
$ws.Range("I4").Value = 0.376631045782902
$ws.Range("J4").Value = 0.3928053077137587
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01499966666666667
$ws.Range("N4").Value = 0.044999
$ws.Range("O4").Value = 0.0008929298633347419
$ws.Range("P4").Value = 0.0009074064137192897
$ws.Range("Q4").Value = 0.02519214016222222
$ws.Range("R4").Value = 0.22672926146
$ws.Range("S4").Value = 0.0003363051082385476
$ws.Range("T4").Value = 0.0003564340555624438
$ws.Range("G5").Value = 1.679513333333333
$ws.Range("H5").Value = 5.038539999999999
$ws.Range("I5").Value = 0.376631045782902
$ws.Range("J5").Value = 0.3928053077137587
$ws.Range("M5").Value = 0.8039865
$ws.Range("N5").Value = 1.607973
$ws.Range("O5").Value = 0.04786130062232345
$ws.Range("P5").Value = 0.03242483195820901
$ws.Range("Q5").Value = 1.35030604657
$ws.Range("R5").Value = 8.101836279419999
$ws.Range("S5").Value = 0.01802605170591554
$ws.Range("T5").Value = 0.01273664609491121
$ws.Range("I6").Value = 0.4998401096732527
$ws.Range("J6").Value = 0.5213055330575571
$ws.Range("M6").Value = 0.1798956666666667
$ws.Range("N6").Value = 0.539687
$ws.Range("O6").Value = 0.01070918551864568
$ws.Range("P6").Value = 0.01088280728907136
$ws.Range("Q6").Value = 0.4009767072252223
$ws.Range("R6").Value = 3.608790365027
$ws.Range("S6").Value = 0.005352880464151066
$ws.Range("T6").Value = 0.005673267654992016
$ws.Range("I7").Value = 0.4998401096732527
$ws.Range("J7").Value = 0.5213055330575571
$ws.Range("N7").Value = 47.39813
$ws.Range("O7").Value = 0.9405365839956962
$ws.Range("P7").Value = 0.9557849543390003
$ws.Range("S7").Value = 0.4701179092961152
$ws.Range("T7").Value = 0.4982559851100855
$ws.Range("I8").Value = 0.4998401096732527
$ws.Range("J8").Value = 0.5213055330575571
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01499966666666667
$ws.Range("N8").Value = 0.044999
$ws.Range("O8").Value = 0.0008929298633347419
$ws.Range("P8").Value = 0.0009074064137192897
$ws.Range("Q8").Value = 0.03343336201988889
$ws.Range("R8").Value = 0.300900258179
$ws.Range("S8").Value = 0.0004463221608197599
$ws.Range("T8").Value = 0.0004730359842037805
$ws.Range("I9").Value = 0.4998401096732527
$ws.Range("J9").Value = 0.5213055330575571
$ws.Range("M9").Value = 0.8039865
$ws.Range("N9").Value = 1.607973
$ws.Range("O9").Value = 0.04786130062232345
$ws.Range("P9").Value = 0.03242483195820901
$ws.Range("Q9").Value = 1.7920379373055
$ws.Range("R9").Value = 10.752227623833
$ws.Range("S9").Value = 0.02392299775216667
$ws.Range("T9").Value = 0.01690324430827586
$ws.Range("G10").Value = 0.550853
$ws.Range("H10").Value = 1.101706
$ws.Range("I10").Value = 0.1235288445438454
$ws.Range("J10").Value = 0.08588915922868416
$ws.Range("M10").Value = 0.1798956666666667
$ws.Range("N10").Value = 0.539687
$ws.Range("O10").Value = 0.01070918551864568
$ws.Range("P10").Value = 0.01088280728907136
$ws.Range("Q10").Value = 0.09909606767033334
$ws.Range("R10").Value = 0.5945764060220001
$ws.Range("S10").Value = 0.001322893313123982
$ws.Range("T10").Value = 0.0009347151681061349
$ws.Range("G11").Value = 0.550853
$ws.Range("H11").Value = 1.101706
$ws.Range("I11").Value = 0.1235288445438454
$ws.Range("J11").Value = 0.08588915922868416
$ws.Range("N11").Value = 47.39813
$ws.Range("O11").Value = 0.9405365839956962
$ws.Range("P11").Value = 0.9557849543390003
$ws.Range("Q11").Value = 8.703134034963334
$ws.Range("R11").Value = 52.21880420978
$ws.Range("S11").Value = 0.1161833974722037
$ws.Range("T11").Value = 0.08209156613160302
$ws.Range("G12").Value = 0.550853
$ws.Range("H12").Value = 1.101706
$ws.Range("I12").Value = 0.1235288445438454
$ws.Range("J12").Value = 0.08588915922868416
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.01499966666666667
$ws.Range("N12").Value = 0.044999
$ws.Range("O12").Value = 0.0008929298633347419
$ws.Range("P12").Value = 0.0009074064137192897
$ws.Range("Q12").Value = 0.008262611382333333
$ws.Range("R12").Value = 0.049575668294
$ws.Range("S12").Value = 0.0001103025942764344
$ws.Range("T12").Value = 0.00007793637395306533
$ws.Range("G13").Value = 0.550853
$ws.Range("H13").Value = 1.101706
$ws.Range("I13").Value = 0.1235288445438454
$ws.Range("J13").Value = 0.08588915922868416
$ws.Range("M13").Value = 0.8039865
$ws.Range("N13").Value = 1.607973
$ws.Range("O13").Value = 0.04786130062232345
$ws.Range("P13").Value = 0.03242483195820901
$ws.Range("Q13").Value = 0.4428783754845
$ws.Range("R13").Value = 1.771513501938
$ws.Range("S13").Value = 0.005912251164241244
$ws.Range("T13").Value = 0.002784941555021941
